# This script updates the "Agroalimentar" domain description and its bullet
# list of topics on the "Eixos" worksheet (row 12, columns C and D), turns on
# wrap text for those two cells, grows the row to fit the new, much longer
# text, and leaves the selection/view on the last edited cell (D12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descricaoAgroalimentar = @'
O domínio Agroalimentar abrange toda a cadeia de produção, transformação e distribuição de alimentos, incluindo tanto os produtos de origem vegetal como animal. Engloba atividades agrícolas, pecuárias, agroindustriais e tecnológicas associadas à alimentação humana e animal.
É um domínio de forte ligação à ciência e à inovação, onde se cruzam áreas como biotecnologia, nutrição, engenharia alimentar, agricultura de precisão e sustentabilidade. Inclui também o desenvolvimento de novos alimentos (funcionais ou personalizados), soluções para conservação, embalagens inteligentes, e tecnologias que garantem a segurança e a rastreabilidade alimentar.
Este domínio procura responder à necessidade de produzir alimentos de forma mais eficiente, saudável, segura e sustentável, reforçando o valor económico e nutricional dos produtos e promovendo a competitividade da fileira agroalimentar portuguesa.
'@

$areasAgroalimentar = @'
	•	Produção agrícola e pecuária
	•	Transformação alimentar e inovação em produtos
	•	Biotecnologia alimentar e nutrição personalizada
	•	Tecnologias de conservação e embalagens inteligentes
	•	Segurança alimentar e rastreabilidade
	•	Sustentabilidade na produção agroindustrial
	•	Valorização de subprodutos e economia circular
	•	Agricultura de precisão e digitalização do setor
'@

# Replace the old "Descrição" (C12) and "Principal área de atuação" (D12)
# text for Agroalimentar with the new, more complete text.
$ws.Range("C12").Value = $descricaoAgroalimentar
$ws.Range("D12").Value = $areasAgroalimentar

# Wrap the text in both cells and enlarge the row so the whole text is
# visible.
$ws.Range("C12:D12").WrapText = $true
$ws.Rows(12).RowHeight = 404

# Match the author's final selection/scroll position on the worksheet.
$ws.Range("D12").Select() | Out-Null
